$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.5979736666666666
$ws.Range("H2").Value = 1.793921
$ws.Range("I2").Value = 0.03342655292740804
$ws.Range("J2").Value = 0.03342655292740804
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 0.03568786572488888
$ws.Range("R2").Value = 0.3211907915239999
$ws.Range("S2").Value = 0.0008700088213402579
$ws.Range("T2").Value = 0.000870008821340258
$ws.Range("G3").Value = 0.5979736666666666
$ws.Range("H3").Value = 1.793921
$ws.Range("I3").Value = 0.03342655292740804
$ws.Range("J3").Value = 0.03342655292740804
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 0.1974642594785556
$ws.Range("R3").Value = 1.777178335307
$ws.Range("S3").Value = 0.004813839218352413
$ws.Range("T3").Value = 0.004813839218352413
$ws.Range("G4").Value = 0.5979736666666666
$ws.Range("H4").Value = 1.793921
$ws.Range("I4").Value = 0.03342655292740804
$ws.Range("J4").Value = 0.03342655292740804
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 1.138009066796333
$ws.Range("R4").Value = 10.242081601167
$ws.Range("S4").Value = 0.02774270488771537
$ws.Range("T4").Value = 0.02774270488771537
$ws.Range("I5").Value = 0.8874158839838097
$ws.Range("J5").Value = 0.8874158839838097
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 0.9474497408848888
$ws.Range("R5").Value = 8.527047667963998
$ws.Range("S5").Value = 0.02309719608061435
$ws.Range("T5").Value = 0.02309719608061436
$ws.Range("I6").Value = 0.8874158839838097
$ws.Range("J6").Value = 0.8874158839838097
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("S6").Value = 0.1277989206541073
$ws.Range("T6").Value = 0.1277989206541073
$ws.Range("I7").Value = 0.8874158839838097
$ws.Range("J7").Value = 0.8874158839838097
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("S7").Value = 0.736519767249088
$ws.Range("T7").Value = 0.7365197672490881
$ws.Range("I8").Value = 0.07915756308878232
$ws.Range("J8").Value = 0.07915756308878232
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05968133333333333
$ws.Range("N8").Value = 0.179044
$ws.Range("O8").Value = 0.02602747651633847
$ws.Range("P8").Value = 0.02602747651633848
$ws.Range("Q8").Value = 0.08451258760533334
$ws.Range("R8").Value = 0.760613288448
$ws.Range("S8").Value = 0.002060271614383863
$ws.Range("T8").Value = 0.002060271614383863
$ws.Range("I9").Value = 0.07915756308878232
$ws.Range("J9").Value = 0.07915756308878232
$ws.Range("O9").Value = 0.144012433133819
$ws.Range("P9").Value = 0.144012433133819
$ws.Range("Q9").Value = 0.4676159582293334
$ws.Range("R9").Value = 4.208543624064
$ws.Range("S9").Value = 0.01139967326135932
$ws.Range("T9").Value = 0.01139967326135932
$ws.Range("I10").Value = 0.07915756308878232
$ws.Range("J10").Value = 0.07915756308878232
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("Q10").Value = 2.694924142976
$ws.Range("S10").Value = 0.06569761821303913
$ws.Range("T10").Value = 0.06569761821303913
